# ISA.xlsx edit
# Commit: "add: se compila a la fpga, se encuentra bug con branch en la fpga, en simulacion sirve"
#
# Content change: the note in E6 ("label = InstrDest - InstrSrc") is replaced
# with a short label "PC Destino", and the fuller explanation is moved into a
# new cell E7 ("Simplemente es la linea de la instr destino") written in a
# smaller (8pt) font, next to the merged "Constante" cell (C7:D7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the old note in E6 with the shorter "PC Destino" label.
# (Keeps E6's existing cell formatting/style untouched.)
$ws.Range("E6").Value2 = "PC Destino"

# Add the longer explanation into the newly-used E7 cell, in a smaller font
# so it fits next to the "Constante" row.
$range_E7 = $ws.Range("E7")
$range_E7.Value2 = "Simplemente es la linea de la instr destino"
$range_E7.Font.Size = 8
$range_E7.Font.ThemeFont = 1  # xlThemeFontMinor -> keeps <scheme val="minor"/> like the rest of the sheet
                               # (leaving Name/Bold/Italic alone keeps them at their Calibri/regular defaults)

# Move the active selection to F7, matching where the author's cursor ended
# up after making the edit.
$ws.Range("F7").Select() | Out-Null
